# Sync file from Google Drive
# Updates EstimatedTimeOfArrival (F), Load (I), OriginCode (K) and
# TypeOfBus (L) values on the NextBus1 / NextBus2 sheets to the refreshed
# Google Drive export.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NextBus1")
$ws.Range("F2").Value = 45684.47931712963
$ws.Range("I2").Value = "SDA"
$ws.Range("F3").Value = 45684.48091435185
$ws.Range("F4").Value = 45684.48255787037
$ws.Range("F5").Value = 45684.4858912037
$ws.Range("F6").Value = 45684.48063657407
$ws.Range("L6").Value = "DD"
$ws.Range("F7").Value = 45684.48263888889
$ws.Range("I7").Value = "SDA"
$ws.Range("F8").Value = 45684.48466435185
$ws.Range("F9").Value = 45684.48015046296
$ws.Range("I9").Value = "SDA"
$ws.Range("F10").Value = 45684.48305555555
$ws.Range("F11").Value = 45684.48041666667
$ws.Range("F12").Value = 45684.481875
$ws.Range("F13").Value = 45684.47857638889
$ws.Range("F14").Value = 45684.48174768518
$ws.Range("K14").Value = 64009
$ws.Range("L14").Value = "SD"
$ws.Range("F15").Value = 45684.48219907407

$ws = $wb.Worksheets.Item("NextBus2")
$ws.Range("F2").Value = 45684.4896875
$ws.Range("I2").Value = "SEA"
$ws.Range("L2").Value = "DD"
$ws.Range("F3").Value = 45684.48743055556
$ws.Range("F4").Value = 45684.49230324074
$ws.Range("F5").Value = 45684.48950231481
$ws.Range("F6").Value = 45684.48627314815
$ws.Range("F7").Value = 45684.49072916667
$ws.Range("I7").Value = "SEA"
$ws.Range("F8").Value = 45684.49001157407
$ws.Range("L8").Value = "SD"
$ws.Range("F9").Value = 45684.48607638889
$ws.Range("I9").Value = "SEA"
$ws.Range("L9").Value = "DD"
$ws.Range("F10").Value = 45684.48752314815
$ws.Range("F11").Value = 45684.49145833333
$ws.Range("L11").Value = "DD"
$ws.Range("F12").Value = 45684.49046296296
$ws.Range("F13").Value = 45684.49385416666
$ws.Range("I13").Value = "SDA"
$ws.Range("L13").Value = "SD"
$ws.Range("F14").Value = 45684.49023148148
$ws.Range("F15").Value = 45684.48811342593
